$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the stray trailing newline character from the body-name cells ---
$ws.Range("A2").Value = "Sun"
$ws.Range("A3").Value = "Mercury"
$ws.Range("A4").Value = "Venus"
$ws.Range("A7").Value = "Jupiter"
$ws.Range("A10").Value = "Neptune"

# --- Column A (body names) becomes Text-formatted ---
$ws.Range("A2:A11").NumberFormat = "@"

# --- Columns B, C, D (except the already explicitly-colored C6) pick up the
#     same "theme 1" font color already used by the header row / column A ---
$ws.Range("B2:B12").Font.ThemeColor = 1
$ws.Range("C2:C5").Font.ThemeColor = 1
$ws.Range("C7:C11").Font.ThemeColor = 1
$ws.Range("D2:D11").Font.ThemeColor = 1
